# error solve ifrs list
# Corrects the financial figures (originally pasted in units of 1 won)
# down to the intended "100 million won" units for 2014-2018 (rows 2-6),
# and removes the erroneous forecast rows for 2019E-2021E (rows 7-9),
# keeping only their identifying columns (A, B, C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1224
$ws.Range("E2").Value = 93
$ws.Range("F2").Value = 93
$ws.Range("G2").Value = 58
$ws.Range("H2").Value = 49
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 2099
$ws.Range("L2").Value = 1298
$ws.Range("M2").Value = 801
$ws.Range("N2").Value = 799
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 414
$ws.Range("Q2").Value = 186
$ws.Range("R2").Value = -92
$ws.Range("S2").Value = -71
$ws.Range("T2").Value = 140
$ws.Range("U2").Value = 47
$ws.Range("V2").Value = 720
$ws.Range("W2").Value = 7.62
$ws.Range("X2").Value = 4.04
$ws.Range("Y2").Value = 7.61
$ws.Range("Z2").Value = 2.34
$ws.Range("AA2").Value = 162.04
$ws.Range("AB2").Value = 86.39
$ws.Range("AC2").Value = 72
$ws.Range("AD2").Value = 14.35
$ws.Range("AE2").Value = 998
$ws.Range("AF2").Value = 1.04
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 82778477

# Row 3
$ws.Range("D3").Value = 1238
$ws.Range("E3").Value = 82
$ws.Range("F3").Value = 82
$ws.Range("G3").Value = 117
$ws.Range("H3").Value = 88
$ws.Range("I3").Value = 88
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2185
$ws.Range("L3").Value = 1150
$ws.Range("M3").Value = 1035
$ws.Range("N3").Value = 1033
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 469
$ws.Range("Q3").Value = 305
$ws.Range("R3").Value = -186
$ws.Range("S3").Value = -131
$ws.Range("T3").Value = 177
$ws.Range("U3").Value = 128
$ws.Range("V3").Value = 436
$ws.Range("W3").Value = 6.59
$ws.Range("X3").Value = 7.13
$ws.Range("Y3").Value = 9.6
$ws.Range("Z3").Value = 4.12
$ws.Range("AA3").Value = 111.15
$ws.Range("AB3").Value = 117.36
$ws.Range("AC3").Value = 101
$ws.Range("AD3").Value = 19.88
$ws.Range("AE3").Value = 1102
$ws.Range("AF3").Value = 1.83
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 93736313

# Row 4
$ws.Range("D4").Value = 1155
$ws.Range("E4").Value = -60
$ws.Range("F4").Value = -60
$ws.Range("G4").Value = -90
$ws.Range("H4").Value = -72
$ws.Range("I4").Value = -72
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2039
$ws.Range("L4").Value = 1162
$ws.Range("M4").Value = 877
$ws.Range("N4").Value = 877
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 479
$ws.Range("Q4").Value = -100
$ws.Range("R4").Value = -152
$ws.Range("S4").Value = 237
$ws.Range("T4").Value = 258
$ws.Range("U4").Value = -359
$ws.Range("V4").Value = 659
$ws.Range("W4").Value = -5.22
$ws.Range("X4").Value = -6.28
$ws.Range("Y4").Value = -7.56
$ws.Range("Z4").Value = -3.43
$ws.Range("AA4").Value = 132.52
$ws.Range("AB4").Value = 79.53
$ws.Range("AC4").Value = -75
$ws.Range("AD4").Value = -17.57
$ws.Range("AE4").Value = 915
$ws.Range("AF4").Value = 1.45
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 95840934

# Row 5
$ws.Range("D5").Value = 996
$ws.Range("E5").Value = 21
$ws.Range("F5").Value = 21
$ws.Range("G5").Value = 9
$ws.Range("H5").Value = 107
$ws.Range("I5").Value = 108
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2178
$ws.Range("L5").Value = 1182
$ws.Range("M5").Value = 996
$ws.Range("N5").Value = 996
$ws.Range("O5").Value = -1
$ws.Range("P5").Value = 483
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = -112
$ws.Range("S5").Value = 101
$ws.Range("T5").Value = 170
$ws.Range("U5").Value = -163
$ws.Range("V5").Value = 755
$ws.Range("W5").Value = 2.08
$ws.Range("X5").Value = 10.78
$ws.Range("Y5").Value = 11.48
$ws.Range("Z5").Value = 5.09
$ws.Range("AA5").Value = 118.7
$ws.Range("AB5").Value = 103.08
$ws.Range("AC5").Value = 112
$ws.Range("AD5").Value = 8.75
$ws.Range("AE5").Value = 1032
$ws.Range("AF5").Value = 0.95
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 96570546

# Row 6
$ws.Range("D6").Value = 820
$ws.Range("E6").Value = -110
$ws.Range("F6").Value = -110
$ws.Range("G6").Value = -231
$ws.Range("H6").Value = -279
$ws.Range("I6").Value = -279
$ws.Range("K6").Value = 2005
$ws.Range("L6").Value = 1311
$ws.Range("M6").Value = 694
$ws.Range("N6").Value = 694
$ws.Range("P6").Value = 483
$ws.Range("Q6").Value = -20
$ws.Range("R6").Value = -60
$ws.Range("S6").Value = 167
$ws.Range("T6").Value = 23
$ws.Range("U6").Value = -43
$ws.Range("V6").Value = 913
$ws.Range("W6").Value = -13.39
$ws.Range("X6").Value = -33.99
$ws.Range("Y6").Value = -32.98
$ws.Range("Z6").Value = -13.33
$ws.Range("AA6").Value = 188.89
$ws.Range("AB6").Value = 45.12
$ws.Range("AC6").Value = -289
$ws.Range("AD6").Value = -6.55
$ws.Range("AE6").Value = 719
$ws.Range("AF6").Value = 2.63
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 96570546

# Rows 7-9: clear all data columns except A, B, C
$ws.Range("D7:AJ9").ClearContents()
